# Weekly update: a new record (week of 2022-08-26, serial 44799) is
# inserted as a new row 88, pushing the existing rows 88..225 down to
# 89..226 (dimension grows from A1:R225 to A1:R226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 88; everything below (old rows
# 88-225) shifts down by one (to 89-226).
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new weekly record.
$ws.Range("A88").Value = 5
$ws.Range("B88").Value = "Macroferia Regional de Talca"
$ws.Range("C88").Value = "Maule"
$ws.Range("D88").Value = 44799
$ws.Range("E88").Value = 7
$ws.Range("F88").Value = 100112017
$ws.Range("G88").Value = "Apio"
$ws.Range("H88").Value = "Americana (o)"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 600
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 10000
$ws.Range("M88").Value = 10000
$ws.Range("N88").Value = "$/docena de matas"
$ws.Range("O88").Value = "Provincia del Elquí"
$ws.Range("P88").Value = 1667
$ws.Range("Q88").Value = 6
$ws.Range("R88").Value = "Hortaliza"
